$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '23.716.13'
$ws.Range("E2").Value = '  +1.29%  '
$ws.Range("D3").Value = '1.653.67'
$ws.Range("E3").Value = '  +1.43%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9998'
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.9996'
$ws.Range("E5").Value = '  -0.14%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '304.12'
$ws.Range("E6").Value = '  +0.24%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3820'
$ws.Range("E7").Value = '  +1.72%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '51.53'
$ws.Range("E8").Value = '  +0.63%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3612'
$ws.Range("E9").Value = '  -0.76%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.251'
$ws.Range("E10").Value = '  +1.58%  '
$ws.Range("E11").Value = '  +0.54%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.0000'
$ws.Range("E12").Value = '  -0.25%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.64'
$ws.Range("E13").Value = '  +1.10%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.549'
$ws.Range("E14").Value = '  +0.18%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.415'
$ws.Range("E15").Value = '  +0.98%  '
$ws.Range("E16").Value = '  -1.08%  '
$ws.Range("D17").Value = '1.653.77'
$ws.Range("E17").Value = '  +1.34%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '97.23'
$ws.Range("E18").Value = '  +3.32%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06975'
$ws.Range("E19").Value = '  +0.06%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.796'
$ws.Range("E20").Value = '  +4.06%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.72'
$ws.Range("E21").Value = '  +0.13%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9995'
$ws.Range("E22").Value = '  -0.14%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.61'
$ws.Range("E23").Value = '  -0.29%  '
$ws.Range("D24").Value = '23.708.62'
$ws.Range("E24").Value = '  +1.28%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.548'
$ws.Range("E25").Value = '  +3.59%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.078'
$ws.Range("E26").Value = '  -1.78%  '
$ws.Range("E27").Value = '  -0.28%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '152.37'
$ws.Range("E28").Value = '  +1.15%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.237'
$ws.Range("E29").Value = '  -1.27%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '135.53'
$ws.Range("E30").Value = '  +1.25%  '
$ws.Range("D31").Value = '1.834.78'
$ws.Range("E31").Value = '  +1.52%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.904'
$ws.Range("E32").Value = '  +1.47%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.087'
$ws.Range("E33").Value = '  +5.32%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '11.94'
$ws.Range("E34").Value = '  +10.33%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.093'
$ws.Range("E35").Value = '  -6.49%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02819'
$ws.Range("E36").Value = '  +1.38%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.2523'
$ws.Range("E37").Value = '  +0.29%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.08841'
$ws.Range("E38").Value = '  +0.73%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.110'
$ws.Range("E39").Value = '  +1.82%  '
$ws.Range("E40").Value = '  -0.76%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '12.83'
$ws.Range("E41").Value = '  +5.59%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.7070'
$ws.Range("E42").Value = '  +0.72%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.341'
$ws.Range("E43").Value = '  -0.37%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '15.96'
$ws.Range("E44").Value = '  -0.49%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6526'
$ws.Range("E45").Value = '  -0.23%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.343'
$ws.Range("E46").Value = '  +2.34%  '
$ws.Range("E47").Value = '  -0.11%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.983'
$ws.Range("E48").Value = '  +0.25%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.07994'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '128.17'
$ws.Range("E50").Value = '  +1.60%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.192'
$ws.Range("E51").Value = '  -0.86%  '
